# Swap the contents of column C (group-code) and column D (group-name)
# for every used row on the active sheet, including the header row.
# This matches the codeforIATI ReportingOrganisationGroup.xlsx update
# which reorders the "code" / "name" shared strings pair-wise, effectively
# swapping the C and D column values throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value()
    $dVal = $dCell.Value()

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
